# Add an "Electrode Locations" column and sort the data rows by electrode
# location (letter prefix, then numeric suffix) - e.g. A1, A4, A6 ... B5,
# C1, C3 ... M5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Read the existing data (file name in col A, value in col B) starting at row 2.
$records = @()
for ($r = 2; $r -le $lastRow; $r++) {
    $fileName = $ws.Cells.Item($r, 1).Value()
    $pmax = $ws.Cells.Item($r, 2).Value()

    if ([string]::IsNullOrEmpty($fileName)) { continue }

    $match = [System.Text.RegularExpressions.Regex]::Match($fileName, '^([A-Za-z]+)(\d+)_')
    $letter = $match.Groups[1].Value
    $number = [int]$match.Groups[2].Value
    $location = "$letter$number"

    $records += [PSCustomObject]@{
        FileName = $fileName
        Pmax     = $pmax
        Location = $location
        Letter   = $letter
        Number   = $number
    }
}

# Sort by electrode letter, then by numeric position (A1-O15 order).
$sorted = $records | Sort-Object Letter, Number

# Header row - give C1 the same (bold/bordered/centered) formatting as B1.
$ws.Cells.Item(1, 3).Value = "Electrode Locations"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

# Write the sorted rows back out.
$row = 2
foreach ($rec in $sorted) {
    $ws.Cells.Item($row, 1).Value = $rec.FileName
    $ws.Cells.Item($row, 2).Value = $rec.Pmax
    $ws.Cells.Item($row, 3).Value = $rec.Location
    $row++
}
